# Update countries & provincias Spain
# - Madagascar's case count overtook "Republica de Africa Central", "Congo"
#   and "Reunion" in the ranking, so those three rows' data shift down by
#   one position and Madagascar's updated stats take the top spot (row 133).
# - The "datos actualizados" timestamp moves from 12:05 to 12:35.
# - Several other countries get refreshed case/recovery/death counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 12:35"

# --- Bielorrusia (row 26) ----------------------------------------------
$ws.Range("B26").Value = 35244
$ws.Range("C26").Value = 941
$ws.Range("D26").Value = 13528
$ws.Range("E26").Value = 21522
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 194

# --- Suiza (row 30) ------------------------------------------------------
$ws.Range("B30").Value = 30725
$ws.Range("C30").Value = 18
$ws.Range("E30").Value = 921
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 1904

# --- Rumania (row 40) ------------------------------------------------------
$ws.Range("B40").Value = 17857
$ws.Range("C40").Value = 145
$ws.Range("D40").Value = 11187
$ws.Range("E40").Value = 5500

# --- Barein (row 54) ------------------------------------------------------
$ws.Range("B54").Value = 8774
$ws.Range("C54").Value = 360
$ws.Range("D54").Value = 4462
$ws.Range("E54").Value = 4300

# --- Marruecos (row 58) ----------------------------------------------------
$ws.Range("B58").Value = 7375
$ws.Range("C58").Value = 43
$ws.Range("D58").Value = 4573
$ws.Range("E58").Value = 2605

# --- Bosnia y Herzegovina (row 83) ------------------------------------------
$ws.Range("B83").Value = 2391
$ws.Range("C83").Value = 19
$ws.Range("D83").Value = 1662
$ws.Range("E83").Value = 588

# --- Rows 133-136: Madagascar jumps ahead of Republica de Africa Central, --
# --- Congo and Reunion; those three now occupy rows 134-136 with their   --
# --- previous values carried down, and Madagascar's fresh stats sit at   --
# --- row 133.                                                            --
$ws.Range("A133").Value = "Madagascar"
$ws.Range("B133").Value = 488
$ws.Range("C133").Value = 40
$ws.Range("D133").Value = 138
$ws.Range("E133").Value = 348
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 2

$ws.Range("A134").Value = "Republica de Africa Central"
$ws.Range("B134").Value = 479
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 18
$ws.Range("E134").Value = 461
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "Congo"
$ws.Range("B135").Value = 469
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 137
$ws.Range("E135").Value = 316
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 16

$ws.Range("A136").Value = "Reunion"
$ws.Range("B136").Value = 449
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 411
$ws.Range("E136").Value = 37
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 1

# --- Cabo Verde (row 141) ---------------------------------------------------
$ws.Range("D141").Value = 130
$ws.Range("E141").Value = 229

Write-Output "done"
